$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 183, pushing the existing rows 183..254 down to 184..255
$row = $ws.Rows.Item(183)
$row.Insert()

# Populate the newly inserted row 183 with the new weekly data point
$ws.Cells.Item(183,1).Value  = 6
$ws.Cells.Item(183,2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(183,3).Value  = "Metropolitana"
$ws.Cells.Item(183,4).Value  = 44489
$ws.Cells.Item(183,5).Value  = 13
$ws.Cells.Item(183,6).Value  = 100112052
$ws.Cells.Item(183,7).Value  = "Albahaca"
$ws.Cells.Item(183,8).Value  = "Sin especificar"
$ws.Cells.Item(183,9).Value  = "Primera"
$ws.Cells.Item(183,10).Value = 150
$ws.Cells.Item(183,11).Value = 5000
$ws.Cells.Item(183,12).Value = 6000
$ws.Cells.Item(183,13).Value = 5467
$ws.Cells.Item(183,14).Value = "$/paquete"
$ws.Cells.Item(183,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(183,16).Value = 5467
$ws.Cells.Item(183,17).Value = 1
$ws.Cells.Item(183,18).Value = "Hortaliza"
